# Ajuste de alguns nomes de campo
# Renames a handful of ER-diagram table-cell field names on slide 1
# (product_id -> productId, age -> birthDate, order_item_id -> orderItemId,
#  order_id -> orderId, product_id -> productId) to match the new
# camelCase convention used elsewhere in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 = "product" entity table -> row 2 ("PK" row), column 2 holds the field name
$productTbl = $s.Shapes.Item(1).Table
$productTbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "productId"

# Shape 3 = "client" entity table -> row 5 holds the "age" field, renamed to "birthDate"
$clientTbl = $s.Shapes.Item(3).Table
$clientTbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "birthDate"

# Shape 5 = "orderItem" entity table -> rows 2-4 hold its PK/FK field names
$orderItemTbl = $s.Shapes.Item(5).Table
$orderItemTbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "orderItemId"
$orderItemTbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "orderId"
$orderItemTbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "productId"
